# Auto-generated edit script applying the scraper refresh described in the diff.
# Updates 'want-to-go' counters (col F / sometimes G) across all four sheets and
# replaces stale listing rows with freshly scraped ones (sheet 展览 rows 13-18,
# sheet 全部类型 row 27).
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1684
$ws.Range("F3").Value = 9610
$ws.Range("F5").Value = 795
$ws.Range("F7").Value = 256
$ws.Range("F10").Value = 79
$ws.Range("F11").Value = 1448
$ws.Range("C13").Value = '上海·NW新界动漫游戏展'
$ws.Range("D13").Value = '长寿路街道澳门路168号 月星家居'
$ws.Range("E13").Value = '2024.06.08 10:00-06.10 16:00'
$ws.Range("F13").Value = 63
$ws.Range("G13").Value = 49
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=83934'
$ws.Range("I13").Value = '//i0.hdslb.com/bfs/openplatform/202404/4eW55lpD1712517166770.jpeg'
$ws.Range("C14").Value = '上海·Redamancy动漫游戏嘉年华-端午篇'
$ws.Range("D14").Value = '中山北路3300号环球港购物中心4楼 上海世嘉都市乐园'
$ws.Range("E14").Value = '2024.06.08 10:00-06.10 17:00'
$ws.Range("F14").Value = 1518
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=84642'
$ws.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202405/Omh4zBq31715321855928.png'
$ws.Range("C15").Value = '上海·cdc动漫展'
$ws.Range("D15").Value = '海潮路133号B1 JUMP工坊'
$ws.Range("F15").Value = 132
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=85110'
$ws.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202405/RMpaP6sF1714725969882.jpeg'
$ws.Range("C16").Value = '上海·夏日欢愉·羽球节庆·原崩铁同人展'
$ws.Range("D16").Value = '鲁班路300号 星光摄影器材城'
$ws.Range("E16").Value = '2024.06.08 10:00-06.09 17:00'
$ws.Range("F16").Value = 321
$ws.Range("G16").Value = 60
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=84742'
$ws.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202404/kbTFe8mJ1713862667234.png'
$ws.Range("C17").Value = '上海·女团驾到·次元女团偶像专区'
$ws.Range("D17").Value = '长宁路1191号来福士西区(W)B1层01号、11号 星零界'
$ws.Range("E17").Value = '2024.06.08 13:00-06.10 18:00'
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 78
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=84796'
$ws.Range("I17").Value = '//i2.hdslb.com/bfs/openplatform/202404/AOS8NlZ31713944402838.jpeg'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '2024-06-08'
$ws.Range("C18").Value = '上海·魔法少女only同人展'
$ws.Range("D18").Value = '南京西路1038号中庭 梅龙镇广场中庭'
$ws.Range("E18").Value = '2024.06.08 09:00-06.08 20:00'
$ws.Range("F18").Value = 161
$ws.Range("G18").Value = 60
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=85738'
$ws.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202405/hlMHqOLK1715830735447.jpeg'
$ws.Range("F19").Value = 420
$ws.Range("F20").Value = 1118
$ws.Range("F21").Value = 108
$ws.Range("F24").Value = 53
$ws.Range("F25").Value = 295
$ws.Range("F30").Value = 650
$ws.Range("F32").Value = 8
$ws.Range("F35").Value = 63
$ws.Range("F36").Value = 193
$ws.Range("F39").Value = 351
$ws.Range("F40").Value = 638
$ws.Range("F43").Value = 336
$ws.Range("F44").Value = 287
$ws.Range("F45").Value = 333
$ws.Range("F46").Value = 66
$ws.Range("F47").Value = 331
$ws.Range("F48").Value = 69

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 12
$ws.Range("F11").Value = 704
$ws.Range("F12").Value = 64
$ws.Range("F16").Value = 54
$ws.Range("F18").Value = 128
$ws.Range("F19").Value = 985
$ws.Range("F21").Value = 1073
$ws.Range("F22").Value = 297
$ws.Range("F23").Value = 669
$ws.Range("F24").Value = 42
$ws.Range("F25").Value = 9
$ws.Range("F27").Value = 334
$ws.Range("F29").Value = 186
$ws.Range("F32").Value = 131
$ws.Range("F39").Value = 38

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 788
$ws.Range("F7").Value = 3803
$ws.Range("F8").Value = 28
$ws.Range("F10").Value = 152
$ws.Range("F11").Value = 134

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1684
$ws.Range("F3").Value = 788
$ws.Range("F4").Value = 9610
$ws.Range("F7").Value = 3803
$ws.Range("F8").Value = 795
$ws.Range("F9").Value = 152
$ws.Range("F10").Value = 152
$ws.Range("F12").Value = 256
$ws.Range("F15").Value = 704
$ws.Range("F16").Value = 1448
$ws.Range("F18").Value = 134
$ws.Range("F19").Value = 1518
$ws.Range("F21").Value = 321
$ws.Range("F22").Value = 161
$ws.Range("F23").Value = 1118
$ws.Range("F24").Value = 108
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = '2024-06-15'
$ws.Range("C27").Value = '上海·THE LAST BLOSSOM~Dreamer''s Band Party 3rd live BLOOM乐队毕业专场'
$ws.Range("D27").Value = '宜昌路179号 万代南梦宫上海文化中心'
$ws.Range("E27").Value = '2024.06.15 15:00-06.15 18:00'
$ws.Range("F27").Value = 54
$ws.Range("G27").Value = 108
$ws.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=86293'
$ws.Range("I27").Value = '//i2.hdslb.com/bfs/openplatform/202405/gG023Oew1716280465865.png'
$ws.Range("F28").Value = 295
$ws.Range("F32").Value = 1073
$ws.Range("F34").Value = 650
$ws.Range("F35").Value = 42
$ws.Range("F36").Value = 9
$ws.Range("F37").Value = 334
$ws.Range("F38").Value = 63
$ws.Range("F41").Value = 352
$ws.Range("F42").Value = 186
$ws.Range("F43").Value = 638
$ws.Range("F46").Value = 336
$ws.Range("F47").Value = 287
$ws.Range("F48").Value = 333
$ws.Range("F50").Value = 38
